$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) titles ---
$ws.Range("C1").Value = "Lectura Inicial"
$ws.Range("D1").Value = "Lectura Final"
$ws.Range("E1").Value = "Causa mes"
$ws.Range("F1").Value = "Consumo mes"
$ws.Range("G1").Value = "Otros"
$ws.Range("H1").Value = "Alumbrado"
$ws.Range("I1").Value = "Kw/h"
$ws.Range("J1").Value = "Valor de Kw/h"
$ws.Range("K1").Value = "Direccion"

# --- Delete the old J column text values (for rows 2-9) before shifting data ---
# First remove the two trailing rows (former rows 8 and 9) entirely.
$ws.Rows("8:9").Delete()

# --- Update remaining data rows (2-7) with new # / Cod restaurante values ---
$ws.Range("A2").Value = 9
$ws.Range("B2").Value = 2

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 2

$ws.Range("A6").Value = 13
$ws.Range("B6").Value = 2

$ws.Range("A7").Value = 14
$ws.Range("B7").Value = 2

# --- Clear the "Causa mes" (old J) text for rows 2 and 3 ---
$ws.Range("J2").ClearContents()
$ws.Range("J3").ClearContents()

# --- Set numeric "Valor de Kw/h" values for rows 4-7 ---
$ws.Range("J4").Value = 602.9728
$ws.Range("J5").Value = 602.9728
$ws.Range("J6").Value = 602.9728
$ws.Range("J7").Value = 602.9728

# --- Set the Direccion value on row 7 ---
$ws.Range("K7").Value = "CRA 13 # 22 B- 11 LC 1-LC 2 LC 00 AV 30 AG"
